$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on cells whose new value is a "clean" numeric-looking
# string (e.g. "1.001", "246.21") so Excel stores it as text instead of auto-converting
# it to a number (which would also mangle things like trailing zeros). Applied per-cell
# since a comma-separated union Range only honors NumberFormat on its first area.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.522.81'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.921.73'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '246.21'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.4845'
$ws.Range("E7").Value = '  +2.63%  '
$ws.Range("D8").Value = '0.2896'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.06709'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '112.11'
$ws.Range("E10").Value = '  +6.46%  '
$ws.Range("D11").Value = '19.33'
$ws.Range("E11").Value = '  +5.08%  '
$ws.Range("D12").Value = '1.922.55'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '0.07589'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '5.365'
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '0.6731'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '295.92'
$ws.Range("E16").Value = '  +1.77%  '
$ws.Range("D17").Value = '30.532.76'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '13.02'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '0.000007550'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.174.31'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.519'
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '6.494'
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("D25").Value = '9.476'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").Value = '164.58'
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("D27").Value = '20.39'
$ws.Range("E27").Value = '  -2.83%  '
$ws.Range("D28").Value = '2.114'
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = '0.1075'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '1.442'
$ws.Range("E30").Value = '  +3.94%  '
$ws.Range("D31").Value = '4.146'
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").Value = '4.054'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").Value = '0.05018'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").Value = '0.7391'
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D36").Value = '0.9999'
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '2.722'
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").Value = '0.02024'
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").Value = '2.696'
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("D40").Value = '110.35'
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("D41").Value = '2.020'
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("D42").Value = '0.4438'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").Value = '0.8662'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("D44").Value = '5.856'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").Value = '70.28'
$ws.Range("E45").Value = '  +4.73%  '
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '7.248'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("D49").Value = '9.232'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").Value = '0.1226'
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").Value = '0.2515'
$ws.Range("E51").Value = '  +2.88%  '
